$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Svornosti" row (row 8) entirely; all rows below shift up.
$ws.Rows.Item(8).Delete()

# Update the view to match: top-left cell A11, active cell A18.
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("A18").Select()
